$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update id column (A) values: add 100 to each existing id
$ws.Range("A2").Value = 101
$ws.Range("A3").Value = 102
$ws.Range("A4").Value = 103
$ws.Range("A5").Value = 104

# Move the selection to B7, matching the author's final cursor position
$ws.Range("B7").Select() | Out-Null
